$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Body Weight Calc")

# Update input values on Sheet1
$ws1.Range("B2").Value = 1.93202989160145          # Battery Weight
$ws1.Range("B3").Value = 0.14435046350608899       # Grip Length
$ws1.Range("B4").Value = 0.34041012843833202       # Bar Length
$ws1.Range("B5").Value = 2.0788034369158401        # Motor Weight

# Update Total Weight formula on Sheet1 (drop B5 term)
$ws1.Range("F2").Formula = "=B2+B6+'Body Weight Calc'!B10+'Body Weight Calc'!B19"

# Update the Housing Height formula on Body Weight Calc sheet (B19)
$ws2.Range("B19").Formula = "=(0.0284*((Sheet1!B4*39.37)^2))-(0.7114*(Sheet1!B4*39.37))+8.2182"

# Update selections to match diff (select Sheet1 last so it remains the active/visible tab)
$ws2.Range("B20").Select()
$ws1.Range("F3").Select()

$wb.Save()
